# Refresh crypto price/volume snapshot (and re-rank a few rows) per
# the GitHub Actions scraper update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''328.31'
$ws.Range("E2").Value = '''0.12%'

# Row 3
$ws.Range("E3").Value = '''1.46%'

# Row 4
$ws.Range("D4").Value = '''5.555'
$ws.Range("E4").Value = '''2.29%'

# Row 5
$ws.Range("D5").Value = '''0.08085'
$ws.Range("E5").Value = '''-0.39%'

# Row 6
$ws.Range("D6").Value = '''1.917'
$ws.Range("E6").Value = '''1.23%'

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '''0.9526'
$ws.Range("E7").Value = '''0.76%'

# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.556'
$ws.Range("E8").Value = '''-7.99%'

# Row 9
$ws.Range("D9").Value = '''0.1183'
$ws.Range("E9").Value = '''0.11%'

# Row 10
$ws.Range("D10").Value = '''0.1852'
$ws.Range("E10").Value = '''-2.06%'

# Row 11
$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").Value = '''10.16'
$ws.Range("E11").Value = '''15.37%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09751'
$ws.Range("E12").Value = '''1.15%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.04488'
$ws.Range("E13").Value = '''6.75%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.1068'
$ws.Range("E14").Value = '''-0.02%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001286'
$ws.Range("E15").Value = '''-1.37%'

# Row 16
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '''0.04201'
$ws.Range("E16").Value = '''-4.31%'

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '''0.005860'
$ws.Range("E17").Value = '''-4.48%'

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.388'
$ws.Range("E18").Value = '''-4.90%'

# Row 19
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").Value = '''4.300'
$ws.Range("E19").Value = '''-0.65%'

# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3462'
$ws.Range("E20").Value = '''-1.57%'

# Row 21
$ws.Range("D21").Value = '''0.1417'
$ws.Range("E21").Value = '''4.05%'

# Row 22
$ws.Range("D22").Value = '''0.2507'
$ws.Range("E22").Value = '''-3.92%'

# Row 23
$ws.Range("D23").Value = '''0.001246'
$ws.Range("E23").Value = '''0.14%'

# Row 24
$ws.Range("D24").Value = '''0.004363'
$ws.Range("E24").Value = '''1.03%'

# Row 25
$ws.Range("D25").Value = '''0.0001191'
$ws.Range("E25").Value = '''-3.98%'

# Row 26
$ws.Range("E26").Value = '''-1.11%'

# Row 38
$ws.Range("D38").Value = '''0.02687'
$ws.Range("E38").Value = '''0.91%'

# Row 39
$ws.Range("D39").Value = '''0.05550'
$ws.Range("E39").Value = '''-0.46%'

# Row 40
$ws.Range("D40").Value = '''0.007574'
$ws.Range("E40").Value = '''-2.69%'

# Row 41
$ws.Range("D41").Value = '''0.1411'
$ws.Range("E41").Value = '''0.73%'

# Row 42
$ws.Range("D42").Value = '''0.008003'
$ws.Range("E42").Value = '''-18.48%'

# Row 43
$ws.Range("D43").Value = '''0.002018'
$ws.Range("E43").Value = '''-5.26%'

# Row 44
$ws.Range("D44").Value = '''0.008408'
$ws.Range("E44").Value = '''-12.63%'

# Row 45
$ws.Range("D45").Value = '''0.00007157'
$ws.Range("E45").Value = '''0.53%'

# Row 46
$ws.Range("D46").Value = '''0.00000000751'
$ws.Range("E46").Value = '''-0.78%'

# Row 47
$ws.Range("D47").Value = '''0.003479'
$ws.Range("E47").Value = '''-0.25%'

# Row 48
$ws.Range("D48").Value = '''0.002272'
$ws.Range("E48").Value = '''-0.77%'

# Row 49
$ws.Range("D49").Value = '''0.00002102'
$ws.Range("E49").Value = '''-0.78%'

# Row 50
$ws.Range("D50").Value = '''0.0002002'
$ws.Range("E50").Value = '''-0.78%'
